# Adds all IG authors as contact
#
# The "Metadata" worksheet lists one "Contact" / "No display for
# ContactDetail" row pair per IG author. Previously there were two such
# rows (rows 10-11); this change adds two more (new rows 12-13), pushing
# everything below (Jurisdiction, Description, Purpose, Copyright, FHIR
# Version, Kind, Type, Base Definition, Abstract, Derivation) down by two
# rows. The "Date" property value is also refreshed to reflect the time
# of this regeneration.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert two new rows right after the existing "Contact" row (row 11),
# shifting all following rows down.
$ws.Rows.Item(12).Resize(2).Insert()

# Populate the two newly inserted "Contact" rows, copying the same
# Property/Value pair used by the existing contact entries.
$ws.Cells.Item(12, 1).Value = "Contact"
$ws.Cells.Item(12, 2).Value = "No display for ContactDetail"
$ws.Cells.Item(12, 1).Style = $ws.Cells.Item(11, 1).Style
$ws.Cells.Item(12, 2).Style = $ws.Cells.Item(11, 2).Style

$ws.Cells.Item(13, 1).Value = "Contact"
$ws.Cells.Item(13, 2).Value = "No display for ContactDetail"
$ws.Cells.Item(13, 1).Style = $ws.Cells.Item(11, 1).Style
$ws.Cells.Item(13, 2).Style = $ws.Cells.Item(11, 2).Style

# Refresh the "Date" metadata property value (column B, row 8).
$ws.Cells.Item(8, 2).Value = "2022-01-21T07:49:24+01:00"
